$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)   # "Diff num of quanv filter"
$ws3 = $wb.Worksheets.Item(3)   # "Diff num of conv filter"

# ---------------------------------------------------------------------------
# Sheet 2 "Diff num of quanv filter"
# Column B rows 9-12 switch from their old fill (orange/red) to the existing
# yellow fill (same style already used elsewhere -> reuse via Interior.Color
# so the engine dedups onto the existing cellXfs entry), and their values
# drop from 2 to 1. Column C rows 10 and 14-19 flip from 1 to 2 (style
# unchanged).
# ---------------------------------------------------------------------------

$ws2.Range("B9").Interior.Color = 65535
$ws2.Range("B9").Value = 1

$ws2.Range("B10").Interior.Color = 65535
$ws2.Range("B10").Value = 1
$ws2.Range("C10").Value = 2

$ws2.Range("B11").Interior.Color = 65535
$ws2.Range("B11").Value = 1

$ws2.Range("B12").Value = 1

$ws2.Range("C14").Value = 2
$ws2.Range("C15").Value = 2
$ws2.Range("C16").Value = 2
$ws2.Range("C17").Value = 2
$ws2.Range("C18").Value = 2
$ws2.Range("C19").Value = 2

# ---------------------------------------------------------------------------
# Sheet 3 "Diff num of conv filter"
# Column C (rows 5-19) gets filled in with values; most rows also pick up a
# green or orange fill + thin border (same look as column B's highlighted
# cells, but without the font flag). Column B rows 11-13 and 15-18 get
# filled in with 2, and B19 switches to the "blue" bordered look used by
# B5/B7/B10. D19 gets a blank white-filled cell.
# ---------------------------------------------------------------------------

# C5 / C6 just need a value - they already carry the right style (s=7).
$ws3.Range("C5").Value = 2
$ws3.Range("C6").Value = 2

# Build the new "green + border, no font" look once, then fan it out with
# PasteSpecial so every cell lands on the very same reused cellXfs entry.
$ws3.Range("C7").Interior.ThemeColor = 10
$ws3.Range("C7").Borders.LineStyle = 1
$ws3.Range("C7").Value = 2
$ws3.Range("C7").Copy()
$ws3.Range("C8").PasteSpecial(-4122)
$ws3.Range("C9").PasteSpecial(-4122)
$ws3.Range("C13").PasteSpecial(-4122)
$ws3.Range("C14").PasteSpecial(-4122)
$ws3.Range("C15").PasteSpecial(-4122)
$ws3.Range("C16").PasteSpecial(-4122)
$ws3.Range("C17").PasteSpecial(-4122)
$ws3.Range("C18").PasteSpecial(-4122)

$ws3.Range("C8").Value = 2
$ws3.Range("C9").Value = 2
$ws3.Range("C13").Value = 2
$ws3.Range("C14").Value = 2
$ws3.Range("C15").Value = 2
$ws3.Range("C16").Value = 1
$ws3.Range("C17").Value = 2
$ws3.Range("C18").Value = 2

# Build the new "orange + border, no font" look once, then fan it out.
$ws3.Range("C10").Interior.ThemeColor = 6
$ws3.Range("C10").Borders.LineStyle = 1
$ws3.Range("C10").Value = 2
$ws3.Range("C10").Copy()
$ws3.Range("C11").PasteSpecial(-4122)
$ws3.Range("C12").PasteSpecial(-4122)
$ws3.Range("C19").PasteSpecial(-4122)

$ws3.Range("C11").Value = 2
$ws3.Range("C12").Value = 2
$ws3.Range("C19").Value = 2

# Column B fills.
$ws3.Range("B11").Value = 2
$ws3.Range("B12").Value = 2
$ws3.Range("B13").Value = 2
$ws3.Range("B15").Value = 2
$ws3.Range("B16").Value = 2
$ws3.Range("B17").Value = 2
$ws3.Range("B18").Value = 2

# B19 switches to the same bordered "blue" look as B5 / B7 / B10.
$ws3.Range("B10").Copy()
$ws3.Range("B19").PasteSpecial(-4122)
$ws3.Range("B19").Value = 2

# D19: new blank cell with a plain white fill (no border, no font flag).
$ws3.Range("D4").Interior.ThemeColor = 2
$ws3.Range("D4").Copy()
$ws3.Range("D19").PasteSpecial(-4122)
$ws3.Range("D4").Clear()

# ---------------------------------------------------------------------------
# Selections - sheet2 ends on E11, sheet3 ends on E10 and stays the active
# (tab-selected) sheet, matching the saved workbook state.
# ---------------------------------------------------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("E11").Select() | Out-Null
$ws3.Activate() | Out-Null
$ws3.Range("E10").Select() | Out-Null
